# Employee DTR report edits (iRipple DTR download) -----------------------
# 1) 05-16-2015 row (row 12): TIME OUT was left blank; fill it in with the
#    same stamp as TIME IN (13:58:27) so the row reads as a single punch.
# 2) 05-18-2015 row (row 14): employee arrived late because of an MRT
#    issue -> charge 1 hour (0.5 day) of Sick Leave, note the remark, and
#    highlight the whole row so it stands out like the other annotated
#    rows in the sheet.
# 3) VL BALANCE legend value updated from 6.4.0 to 5.4.0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 12 (05-16-2015): copy TIME IN into the empty TIME OUT cell ---
$ws.Range("D12").Value = $ws.Range("C12").Value2

# --- 2. Row 14 (05-18-2015): SL charge, remark, and highlight ---
$ws.Range("A14:P14").Interior.Color = 16753510
$ws.Range("I14").Value = 0.5
$ws.Range("P14").Value = "~SL 1 hour late due to MRT technical problem. (Please charge to VL.)"

# --- 3. VL balance figure ---
$ws.Range("C30").Value = "5.4.0"
